$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2025-04-04 Friday" "2025-04-05 Saturday"

Replace-Text "519×7=" "773×6="
Replace-Text "168×8=" "417×7="
Replace-Text "256×3=" "688×7="
Replace-Text "545×5=" "106×9="
Replace-Text "760×3=" "571×9="
Replace-Text "973×4=" "213×6="
Replace-Text "631×6=" "648×2="
Replace-Text "334×7=" "561×2="
Replace-Text "286×6=" "246×3="
Replace-Text "320×9=" "931×2="
Replace-Text "108×2=" "615×2="
Replace-Text "473×9=" "628×8="
Replace-Text "581×6=" "237×9="
Replace-Text "850×8=" "864×9="
Replace-Text "450×3=" "584×7="
Replace-Text "630×4=" "410×3="
Replace-Text "138×9=" "128×3="
Replace-Text "268×5=" "729×4="
Replace-Text "684×4=" "521×2="
Replace-Text "406×9=" "396×9="
Replace-Text "432×2=" "518×6="
Replace-Text "187×3=" "617×2="
Replace-Text "142×5=" "438×5="
Replace-Text "669×6=" "718×2="
Replace-Text "365×5=" "797×7="
